$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "Kichwa cha Video" "Video Title"
Replace-Text "Mada" "Topic"
Replace-Text "Malengo" "Aim(s)"
Replace-Text "Urefu" "Length"
Replace-Text "Mahali pa Kambi" "Camp Location"
Replace-Text "Wawezeshaji" "Facilitators"
Replace-Text "N. ya wanafunzi" "N. of students"
Replace-Text "Tarehe" "Date"
Replace-Text "Rasilimali" "Resources"
Replace-Text "inahitajika" "needed"
Replace-Text "Maandalizi" "Preparations"
Replace-Text "Muda wa video" "Video time"
Replace-Text "Mwezeshaji anafanya nini" "What facilitator does"
Replace-Text "Wanachofanya wanafunzi" "What learners do"
Replace-Text "Utangulizi Mkuu wa Video ya VMC" "General VMC Video Introduction"
Replace-Text "Utangulizi wa Video" "Video Introduction"
Replace-Text "Utangulizi wa jaribio la kwanza" "Introduction of the first experiment"
Replace-Text "Kusaidia mchakato, kuchochea mawazo" "Assist the process, provoke thoughts"

# Document default language: Swahili (Kenya) -> Swahili (Tanzania).
# The runtime's COM surface has no property that reaches
# styles.xml's w:docDefaults/w:rPrDefault directly, so we set the
# language on the base "Normal" style (w:default="1", no basedOn),
# which is the documented Style.LanguageID property and - since every
# paragraph in this document implicitly uses Normal - overrides the
# inherited docDefaults language for the whole document, matching the
# intended sw-KE -> sw-TZ change.
$d.Styles("Normal").LanguageID = "sw-TZ"
